# error solve ifrs list
# Replaces the financial figures for 우진플라임 (rows 2-6) with corrected
# values, and removes the erroneous estimate columns/rows (AG/AH in rows
# 4-6, and the whole D:AI block in rows 7-9) that were based on bad data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 1898
$ws.Range("E2").Value = -72
$ws.Range("F2").Value = -72
$ws.Range("G2").Value = -34
$ws.Range("H2").Value = -2
$ws.Range("I2").Value = -2
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 3575
$ws.Range("L2").Value = 2747
$ws.Range("M2").Value = 827
$ws.Range("N2").Value = 827
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 50
$ws.Range("Q2").Value = -25
$ws.Range("R2").Value = -1310
$ws.Range("S2").Value = 1177
$ws.Range("T2").Value = 1374
$ws.Range("U2").Value = -1399
$ws.Range("V2").Value = 1971
$ws.Range("W2").Value = -3.79
$ws.Range("X2").Value = -0.09
$ws.Range("Y2").Value = -0.21
$ws.Range("Z2").Value = -0.06
$ws.Range("AA2").Value = 332.02
$ws.Range("AB2").Value = 1494.96
$ws.Range("AC2").Value = -17
$ws.Range("AD2").Value = -241.12
$ws.Range("AE2").Value = 8274
$ws.Range("AF2").Value = 0.5
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 10000000

# --- Row 3 ---
$ws.Range("D3").Value = 1832
$ws.Range("E3").Value = -107
$ws.Range("F3").Value = -107
$ws.Range("G3").Value = -210
$ws.Range("H3").Value = -173
$ws.Range("I3").Value = -173
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2943
$ws.Range("L3").Value = 2286
$ws.Range("M3").Value = 658
$ws.Range("N3").Value = 658
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 50
$ws.Range("Q3").Value = -63
$ws.Range("R3").Value = 364
$ws.Range("S3").Value = -285
$ws.Range("T3").Value = 172
$ws.Range("U3").Value = -235
$ws.Range("V3").Value = 1687
$ws.Range("W3").Value = -5.85
$ws.Range("X3").Value = -9.42
$ws.Range("Y3").Value = -23.23
$ws.Range("Z3").Value = -5.29
$ws.Range("AA3").Value = 347.48
$ws.Range("AB3").Value = 1159.1
$ws.Range("AC3").Value = -1725
$ws.Range("AD3").Value = -2.51
$ws.Range("AE3").Value = 6578
$ws.Range("AF3").Value = 0.66
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 10000000

# --- Row 4 (AG4/AH4 removed entirely) ---
$ws.Range("D4").Value = 2283
$ws.Range("E4").Value = 161
$ws.Range("F4").Value = 161
$ws.Range("G4").Value = 157
$ws.Range("H4").Value = 148
$ws.Range("I4").Value = 148
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3040
$ws.Range("L4").Value = 2241
$ws.Range("M4").Value = 799
$ws.Range("N4").Value = 799
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 50
$ws.Range("Q4").Value = 285
$ws.Range("R4").Value = 27
$ws.Range("S4").Value = -184
$ws.Range("T4").Value = 56
$ws.Range("U4").Value = 229
$ws.Range("V4").Value = 1506
$ws.Range("W4").Value = 7.07
$ws.Range("X4").Value = 6.48
$ws.Range("Y4").Value = 20.32
$ws.Range("Z4").Value = 4.95
$ws.Range("AA4").Value = 280.45
$ws.Range("AB4").Value = 1449.83
$ws.Range("AC4").Value = 1480
$ws.Range("AD4").Value = 5.38
$ws.Range("AE4").Value = 7991
$ws.Range("AF4").Value = 1
$ws.Range("AG4").ClearContents()
$ws.Range("AH4").ClearContents()
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 10000000

# --- Row 5 (AG5/AH5 removed entirely) ---
$ws.Range("D5").Value = 2296
$ws.Range("E5").Value = 62
$ws.Range("F5").Value = 62
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3131
$ws.Range("L5").Value = 2320
$ws.Range("M5").Value = 811
$ws.Range("N5").Value = 811
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 50
$ws.Range("Q5").Value = -134
$ws.Range("R5").Value = -106
$ws.Range("S5").Value = 107
$ws.Range("T5").Value = 180
$ws.Range("U5").Value = -315
$ws.Range("V5").Value = 1612
$ws.Range("W5").Value = 2.71
$ws.Range("X5").Value = 0.32
$ws.Range("Y5").Value = 0.92
$ws.Range("Z5").Value = 0.24
$ws.Range("AA5").Value = 285.93
$ws.Range("AB5").Value = 1471.37
$ws.Range("AC5").Value = 74
$ws.Range("AD5").Value = 115.71
$ws.Range("AE5").Value = 8114
$ws.Range("AF5").Value = 1.06
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 10000000

# --- Row 6 (no J6/O6 before or after; AG6/AH6 removed entirely) ---
$ws.Range("D6").Value = 1719
$ws.Range("E6").Value = -118
$ws.Range("F6").Value = -118
$ws.Range("G6").Value = -145
$ws.Range("H6").Value = -136
$ws.Range("I6").Value = -136
$ws.Range("K6").Value = 2877
$ws.Range("L6").Value = 2210
$ws.Range("M6").Value = 667
$ws.Range("N6").Value = 667
$ws.Range("P6").Value = 50
$ws.Range("Q6").Value = 123
$ws.Range("R6").Value = -67
$ws.Range("S6").Value = -67
$ws.Range("T6").Value = 77
$ws.Range("U6").Value = 46
$ws.Range("V6").Value = 1545
$ws.Range("W6").Value = -6.85
$ws.Range("X6").Value = -7.91
$ws.Range("Y6").Value = -18.39
$ws.Range("Z6").Value = -4.52
$ws.Range("AA6").Value = 331.16
$ws.Range("AB6").Value = 1185.98
$ws.Range("AC6").Value = -1359
$ws.Range("AD6").Value = -3.84
$ws.Range("AE6").Value = 6673
$ws.Range("AF6").Value = 0.78
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 10000000

# --- Rows 7, 8, 9: drop all the estimate data, keep only A/B/C ---
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
